# "remove column from alcohol data"
#
# The data sheet (Sheet1) had an extra column (M) of alcohol measurements
# that duplicated/shifted the following column (N). The fix removes column
# M entirely, so the old column N (and its values) slides left to become
# the new column M. This shrinks the used range from A1:N119 to A1:M119.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete column M - remaining columns to the right (N) shift one place left.
$ws.Columns("M:M").Delete()

# Leave the selection where the deleted column used to be, now occupied by
# the shifted data (matches the post-edit active cell in the workbook).
$ws.Range("M1").Select()
